$wb = $excel.ActiveWorkbook

# --- Sheet "Masters" ---
$ws2 = $wb.Worksheets.Item("Masters")

$ws2.Range("H2").Value = 0
$ws2.Range("H3").ClearContents()

$ws2.Range("H3").Select()

# --- Sheet "muži" (results of round 1 "1.kolo" get updated, round 2 "2.kolo" partly) ---
$ws1 = $wb.Worksheets.Item("muži")

$ws1.Range("E2").Value = 2.2999999999999998
$ws1.Range("F2").Value = 1.2

$ws1.Range("E3").Value = 5.6
$ws1.Range("F3").Value = 0

$ws1.Range("E4:E18").Value = 12.5

$ws1.Activate()
$ws1.Range("F4").Select()
